$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: update Objetivos B/C text (was misaligned with Robson professor name)
$ws.Range("B10").Value = 'Fornecer ao aluno os critérios básicos para elaboração das estações de tratamento de água utilizando diferentes tecnologias de tratamento de água para consumo humano e dos resíduos gerados. Os estudantess irão elaborar projetos utilizando as diversas tecnologias de tratamento.'
$ws.Range("C10").Value = 'Fornecer ao aluno os critérios básicos para elaboração das estações de tratamento de água utilizando diferentes tecnologias de tratamento de água para consumo humano e dos resíduos gerados. Os estudantess irão elaborar projetos utilizando as diversas tecnologias de tratamento.'

# Clear rows 13-21 completely before rebuilding the shifted table
$ws.Range("A13:C21").Clear()

# Row 13
$ws.Range("B13").Value = '7455355 - Robson da Silva Rocha'
$ws.Range("B13").Font.Bold = $false
$ws.Range("B13").WrapText = $true
$ws.Range("B13").VerticalAlignment = -4160
$ws.Range("C13").Value = '7455355 - Robson da Silva Rocha'
$ws.Range("C13").Font.Bold = $false
$ws.Range("C13").WrapText = $true
$ws.Range("C13").VerticalAlignment = -4160
$ws.Range("C13").Font.Color = 255
$ws.Rows.Item(13).AutoFit()

# Row 14
$ws.Range("A14").Value = 'Programa resumido:'
$ws.Range("A14").Font.Bold = $true
$ws.Range("A14").WrapText = $false
$ws.Range("A14").VerticalAlignment = -4160
$ws.Range("B14").Value = 'Concepção de Sistemas de Tratamento de Água em Função da Qualidade da Água Bruta; Projeto de ETAs de Ciclo Completo com Emprego da Decantação ou da Flotação por ar Dissolvido por Clarificação; Projeto de ETAs de Filtração Direta Descendente; Projeto de ETAs de Filtração Direta Ascendente; Projeto de ETAs de Dupla Filtração; Projeto de ETAs por Floto-Filtração; Projeto de ETAs de filtração em Múltiplas Etapas - FiME; Métodos Alternativos de Desinfecção e Adsorção em Carvão Ativado; Tratamento dos Resíduos Gerados nas ETAs e Reuso da Água Recuperada.'
$ws.Range("B14").Font.Bold = $false
$ws.Range("B14").WrapText = $true
$ws.Range("B14").VerticalAlignment = -4160
$ws.Range("C14").Value = 'Concepção de Sistemas de Tratamento de Água em Função da Qualidade da Água Bruta; Projeto de ETAs de Ciclo Completo com Emprego da Decantação ou da Flotação por ar Dissolvido por Clarificação; Projeto de ETAs de Filtração Direta Descendente; Projeto de ETAs de Filtração Direta Ascendente; Projeto de ETAs de Dupla Filtração; Projeto de ETAs por Floto-Filtração; Projeto de ETAs de filtração em Múltiplas Etapas - FiME; Métodos Alternativos de Desinfecção e Adsorção em Carvão Ativado; Tratamento dos Resíduos Gerados nas ETAs e Reuso da Água Recuperada.'
$ws.Range("C14").Font.Bold = $false
$ws.Range("C14").WrapText = $true
$ws.Range("C14").VerticalAlignment = -4160
$ws.Range("C14").Font.Color = 255
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Range("A15").Value = 'Short syllabus:'
$ws.Range("A15").Font.Bold = $true
$ws.Range("A15").WrapText = $false
$ws.Range("A15").VerticalAlignment = -4160
$ws.Rows.Item(15).RowHeight = 60

# Row 16
$ws.Range("A16").Value = 'Programa:'
$ws.Range("A16").Font.Bold = $true
$ws.Range("A16").WrapText = $false
$ws.Range("A16").VerticalAlignment = -4160
$ws.Range("B16").Value = '1 - Concepção de Sistemas de Tratamento de Água em Função da Qualidade da Água Bruta - Tecnologias de Tratamento de Água e dos Resíduos Gerados nos ETAs; 2 - Projeto de ETAs de Ciclo Completo com Emprego da Decantação ou da Flotação por ar Dissolvido para Clarificação; Características de água bruta; características de coagulação e coagulantes; 3 - Projeto de ETAs de Filtração Direta Descendente; características de água bruta; características da coagulação e da floculação; efeito da floculação; filtros com taxa constante e taxa declinante; mecanismo da coagulação e principais coagulantes; 4 - Projeto de ETAs de Filtração Direta Ascendente; características de água bruta; mecanismo da coagulação e principais coagulantes; 5 - Projeto de ETAs de Dupla Filtração; Características de água bruta; mecanismo da coagulação e principais coagulantes; instalação com baterias independentes de filtros ascendentes e descendentes; instalação com filtros ascendentes/descendentes;; 6 - Projeto de ETAs por Floto-Filtração; características de água bruta; características da coagulação e da floculação; características dos filtros; 7 - Projeto de ETAs de Filtração em Múltiplas Etapas - FiME; características de água bruta; instalações de pré-filtração dinâmica; pré-filtração em pedregulho com escoamento ascendente, descendente ou horizontal e filtração lenta em areia; considerações sobre a operação e manutenção; 8 - Métodos Alternativos de Desinfecção e Adsorção em Carvão Ativado; unidades de pré e de pós-desinfecção; características da água e formação sub-produtos; isotermas de adsorção; parâmetros de projeto de adosrção e da câmara de contato; 9 - Tratamento dos resíduos Gerados na ETAs e Reuso da Água Recuperada; tecnologia de tratamento de água e características do sistema de tratamento dos resíduos.'
$ws.Range("B16").Font.Bold = $false
$ws.Range("B16").WrapText = $true
$ws.Range("B16").VerticalAlignment = -4160
$ws.Range("C16").Value = '1 - Concepção de Sistemas de Tratamento de Água em Função da Qualidade da Água Bruta - Tecnologias de Tratamento de Água e dos Resíduos Gerados nos ETAs; 2 - Projeto de ETAs de Ciclo Completo com Emprego da Decantação ou da Flotação por ar Dissolvido para Clarificação; Características de água bruta; características de coagulação e coagulantes; 3 - Projeto de ETAs de Filtração Direta Descendente; características de água bruta; características da coagulação e da floculação; efeito da floculação; filtros com taxa constante e taxa declinante; mecanismo da coagulação e principais coagulantes; 4 - Projeto de ETAs de Filtração Direta Ascendente; características de água bruta; mecanismo da coagulação e principais coagulantes; 5 - Projeto de ETAs de Dupla Filtração; Características de água bruta; mecanismo da coagulação e principais coagulantes; instalação com baterias independentes de filtros ascendentes e descendentes; instalação com filtros ascendentes/descendentes;; 6 - Projeto de ETAs por Floto-Filtração; características de água bruta; características da coagulação e da floculação; características dos filtros; 7 - Projeto de ETAs de Filtração em Múltiplas Etapas - FiME; características de água bruta; instalações de pré-filtração dinâmica; pré-filtração em pedregulho com escoamento ascendente, descendente ou horizontal e filtração lenta em areia; considerações sobre a operação e manutenção; 8 - Métodos Alternativos de Desinfecção e Adsorção em Carvão Ativado; unidades de pré e de pós-desinfecção; características da água e formação sub-produtos; isotermas de adsorção; parâmetros de projeto de adosrção e da câmara de contato; 9 - Tratamento dos resíduos Gerados na ETAs e Reuso da Água Recuperada; tecnologia de tratamento de água e características do sistema de tratamento dos resíduos.'
$ws.Range("C16").Font.Bold = $false
$ws.Range("C16").WrapText = $true
$ws.Range("C16").VerticalAlignment = -4160
$ws.Range("C16").Font.Color = 255
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Range("A17").Value = 'Syllabus:'
$ws.Range("A17").Font.Bold = $true
$ws.Range("A17").WrapText = $false
$ws.Range("A17").VerticalAlignment = -4160
$ws.Rows.Item(17).RowHeight = 120

# Row 18
$ws.Range("A18").Value = 'Avaliação:'
$ws.Range("A18").Font.Bold = $true
$ws.Range("A18").WrapText = $false
$ws.Range("A18").VerticalAlignment = -4160
$ws.Rows.Item(18).AutoFit()

# Row 19
$ws.Range("A19").Value = 'Método:'
$ws.Range("A19").Font.Bold = $true
$ws.Range("A19").WrapText = $false
$ws.Range("A19").VerticalAlignment = -4160
$ws.Range("B19").Value = 'Aulas expositivas, estudos de projetos sobre as diferentes tecnologias de tratamento; visitas técnicas.'
$ws.Range("B19").Font.Bold = $false
$ws.Range("B19").WrapText = $true
$ws.Range("B19").VerticalAlignment = -4160
$ws.Range("C19").Value = 'Aulas expositivas, estudos de projetos sobre as diferentes tecnologias de tratamento; visitas técnicas.'
$ws.Range("C19").Font.Bold = $false
$ws.Range("C19").WrapText = $true
$ws.Range("C19").VerticalAlignment = -4160
$ws.Range("C19").Font.Color = 255
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Range("A20").Value = 'Critério:'
$ws.Range("A20").Font.Bold = $true
$ws.Range("A20").WrapText = $false
$ws.Range("A20").VerticalAlignment = -4160
$ws.Range("B20").Value = 'Avaliação composta por 3 (três) provas, sendo uma substitutiva, e por exercícios sobre as unidades de uma estação de tratamento de água.
Nota Final = 0,4 x MP + 0,6 x MT
MP: média das provas; ME: média de trabalhos 
* valor mínimo da média das notas das provas (MP) = 5,0
* valor mínimo da média das notas dos trabalhos e projetos (MT) = 5,0'
$ws.Range("B20").Font.Bold = $false
$ws.Range("B20").WrapText = $true
$ws.Range("B20").VerticalAlignment = -4160
$ws.Range("C20").Value = 'Avaliação composta por 3 (três) provas, sendo uma substitutiva, e por exercícios sobre as unidades de uma estação de tratamento de água.
Nota Final = 0,4 x MP + 0,6 x MT
MP: média das provas; ME: média de trabalhos 
* valor mínimo da média das notas das provas (MP) = 5,0
* valor mínimo da média das notas dos trabalhos e projetos (MT) = 5,0'
$ws.Range("C20").Font.Bold = $false
$ws.Range("C20").WrapText = $true
$ws.Range("C20").VerticalAlignment = -4160
$ws.Range("C20").Font.Color = 255
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Range("A21").Value = 'Norma de recuperação:'
$ws.Range("A21").Font.Bold = $true
$ws.Range("A21").WrapText = $false
$ws.Range("A21").VerticalAlignment = -4160
$ws.Range("B21").Value = 'Prova única com nota igual ou superior a 5,0.'
$ws.Range("B21").Font.Bold = $false
$ws.Range("B21").WrapText = $true
$ws.Range("B21").VerticalAlignment = -4160
$ws.Range("C21").Value = 'Prova única com nota igual ou superior a 5,0.'
$ws.Range("C21").Font.Bold = $false
$ws.Range("C21").WrapText = $true
$ws.Range("C21").VerticalAlignment = -4160
$ws.Range("C21").Font.Color = 255
$ws.Rows.Item(21).RowHeight = 60

# Row 22
$ws.Range("A22").Value = 'Bibliografia:'
$ws.Range("A22").Font.Bold = $true
$ws.Range("A22").WrapText = $false
$ws.Range("A22").VerticalAlignment = -4160
$ws.Range("B22").Value = 'DI BERNARDO, L. Métodos e Técnicas de Tratamento de Água. ASSOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL & LUIZ DI BERNARDO 2 V., Rio de Janeiro, 1993 (2005)
DI BERNARDO, L. Algas e suas Influências na Qualidade da Água e nas Tecnologias de Tratamento ASSOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL & LUIZ DI BERNARDO, Rio de Janeiro, 1995.
PROGRAMA DE PESQUISA EM SANEAMENTO BÁSICO Tratamento de Água de Abastecimento por Filtração em Múltiplas Etapas ASSOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 1999 (PROGRAMA DE PESQUISA EM SANEAMENTO BÁSICO. Noções Gerais de Tratamento e Disposição Final de Lodos de Estações de Tratamento de Água ASOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 2000.'
$ws.Range("B22").Font.Bold = $false
$ws.Range("B22").WrapText = $true
$ws.Range("B22").VerticalAlignment = -4160
$ws.Range("C22").Value = 'DI BERNARDO, L. Métodos e Técnicas de Tratamento de Água. ASSOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL & LUIZ DI BERNARDO 2 V., Rio de Janeiro, 1993 (2005)
DI BERNARDO, L. Algas e suas Influências na Qualidade da Água e nas Tecnologias de Tratamento ASSOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL & LUIZ DI BERNARDO, Rio de Janeiro, 1995.
PROGRAMA DE PESQUISA EM SANEAMENTO BÁSICO Tratamento de Água de Abastecimento por Filtração em Múltiplas Etapas ASSOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 1999 (PROGRAMA DE PESQUISA EM SANEAMENTO BÁSICO. Noções Gerais de Tratamento e Disposição Final de Lodos de Estações de Tratamento de Água ASOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 2000.'
$ws.Range("C22").Font.Bold = $false
$ws.Range("C22").WrapText = $true
$ws.Range("C22").VerticalAlignment = -4160
$ws.Range("C22").Font.Color = 255
$ws.Rows.Item(22).RowHeight = 120
